$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the D-column cells whose new price text is numeric-looking
# (e.g. "1.006", "0.06820") to stay stored as literal text, matching
# the workbook authoring convention where every Price cell is an
# inline/shared text string rather than a real number.
$textRows = @(4, 5, 6, 7, 8, 9, 10, 11, 13, 14, 15, 16, 18, 19, 22, 23, 24, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49, 50, 51)
foreach ($r in $textRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "29.514.86"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "1.900.40"
$ws.Range("E3").Value = "  -0.92%  "
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "339.09"
$ws.Range("E5").Value = "  +4.18%  "
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("D7").Value = "0.4753"
$ws.Range("E7").Value = "  -1.45%  "
$ws.Range("D8").Value = "0.4005"
$ws.Range("E8").Value = "  -1.56%  "
$ws.Range("D9").Value = "0.08048"
$ws.Range("E9").Value = "  -2.04%  "
$ws.Range("D10").Value = "0.9913"
$ws.Range("E10").Value = "  -2.18%  "
$ws.Range("D11").Value = "23.22"
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("D12").Value = "1.905.69"
$ws.Range("E12").Value = "  -2.07%  "
$ws.Range("D13").Value = "5.946"
$ws.Range("E13").Value = "  -1.99%  "
$ws.Range("D14").Value = "7.102"
$ws.Range("E14").Value = "  -2.04%  "
$ws.Range("D15").Value = "89.11"
$ws.Range("E15").Value = "  -2.88%  "
$ws.Range("D16").Value = "0.06820"
$ws.Range("E16").Value = "  -0.76%  "
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").Value = "0.00001020"
$ws.Range("E18").Value = "  -1.93%  "
$ws.Range("D19").Value = "17.34"
$ws.Range("E19").Value = "  -1.59%  "
$ws.Range("E20").Value = "  -0.45%  "
$ws.Range("D21").Value = "29.525.94"
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").Value = "5.511"
$ws.Range("E22").Value = "  -2.77%  "
$ws.Range("D23").Value = "11.62"
$ws.Range("E23").Value = "  -1.00%  "
$ws.Range("D24").Value = "2.153"
$ws.Range("E24").Value = "  -1.59%  "
$ws.Range("D25").Value = "2.155.45"
$ws.Range("E25").Value = "  -0.73%  "
$ws.Range("D26").Value = "157.08"
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("D27").Value = "6.509"
$ws.Range("E27").Value = "  -2.40%  "
$ws.Range("D28").Value = "19.66"
$ws.Range("E28").Value = "  -1.77%  "
$ws.Range("D29").Value = "2.056"
$ws.Range("E29").Value = "  -2.82%  "
$ws.Range("D30").Value = "119.16"
$ws.Range("E30").Value = "  -1.46%  "
$ws.Range("D31").Value = "0.9967"
$ws.Range("E31").Value = "  -1.80%  "
$ws.Range("D32").Value = "0.09536"
$ws.Range("E32").Value = "  -0.86%  "
$ws.Range("D33").Value = "5.482"
$ws.Range("E33").Value = "  -2.92%  "
$ws.Range("D34").Value = "1.389"
$ws.Range("E34").Value = "  +1.09%  "
$ws.Range("D35").Value = "3.531"
$ws.Range("E36").Value = "  +4.57%  "
$ws.Range("D37").Value = "0.02246"
$ws.Range("E37").Value = "  -1.66%  "
$ws.Range("D38").Value = "1.201"
$ws.Range("E38").Value = "  +1.72%  "
$ws.Range("D39").Value = "0.5824"
$ws.Range("E39").Value = "  -2.47%  "
$ws.Range("D40").Value = "10.56"
$ws.Range("E40").Value = "  -3.19%  "
$ws.Range("D41").Value = "7.735"
$ws.Range("E41").Value = "  -4.12%  "
$ws.Range("D42").Value = "0.1819"
$ws.Range("E42").Value = "  -1.52%  "
$ws.Range("D43").Value = "2.434"
$ws.Range("E43").Value = "  +1.28%  "
$ws.Range("D44").Value = "1.268"
$ws.Range("E44").Value = "  -1.04%  "
$ws.Range("D45").Value = "12.12"
$ws.Range("E45").Value = "  -2.66%  "
$ws.Range("D46").Value = "0.07365"
$ws.Range("E46").Value = "  -3.06%  "
$ws.Range("D47").Value = "0.5495"
$ws.Range("E47").Value = "  -1.78%  "
$ws.Range("D48").Value = "1.954"
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("D49").Value = "116.32"
$ws.Range("E49").Value = "  -2.03%  "
$ws.Range("D50").Value = "2.376"
$ws.Range("E50").Value = "  -2.04%  "
$ws.Range("D51").Value = "71.15"
$ws.Range("E51").Value = "  -1.58%  "
